$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove all existing hyperlinks (and their relationships) so we can
#     rebuild the hyperlink set cleanly for the new data. ---
$ws.Hyperlinks.Delete()

# --- Drop the last two login rows (rows 5 and 6) ---
$ws.Range("A5:B6").EntireRow.Delete()

# --- Overwrite the remaining data rows with the new login records.
#     Order of assignment controls the order new shared strings are
#     interned in, so values are set row/cell in the sequence that
#     reproduces the target shared-string table. ---
$ws.Range("A3").Value = "mukraheel@gmail.com"
$ws.Range("B3").Value = "R@heel123"
$ws.Range("B2").Value = "Raheel1234"
$ws.Range("A4").Value = "Testmail@gmail.com"
$ws.Range("B4").Value = "test123456"
$ws.Range("A2").Value = "accraheel123@gmail.com"

# --- Re-create the hyperlinks against the new email / password cells,
#     in the order that yields rId1..rId4. ---
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:mukraheel@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:R@heel123")
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:accraheel123@gmail.com")
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:Testmail@gmail.com")

# Hyperlinks.Add stamps a freshly-minted "Hyperlink" cell style on every
# cell it touches; collapse them all back onto the single pre-existing
# "Hyperlink" cell style (style index 1) that the sheet already used.
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"
$ws.Range("A4").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"

# --- New column B width (about 11 characters) ---
$ws.Columns("B").ColumnWidth = 10.14

# --- Selection moves to O8 ---
$ws.Range("O8").Select()
